$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.658.97'

$ws.Range("D3").Value = '2.277.07'
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.66'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  +8.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.13'
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = '  -1.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = '  +1.46%  '

$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = '  +1.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.48'
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = '  +5.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = '  +0.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.86'
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = '  +11.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.64'
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = '  +0.54%  '

$ws.Range("D15").Value = '2.621.75'
$ws.Range("E15").Value = '  -0.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.873'
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = '  +2.25%  '

$ws.Range("D17").Value = '2.279.34'
$ws.Range("E17").Value = '  -0.40%  '

$ws.Range("D18").Value = '43.577.17'
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("E19").Value = '  -1.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.97'
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = '  +11.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.02'
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = '  -0.30%  '

$ws.Range("E22").Value = '  -3.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.98'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  +9.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '232.25'
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = '  -0.35%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.55'
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = '  +3.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.61'
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = '  +2.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.40'
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("E30").Value = '  -1.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.99'
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = '  -2.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.46'
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = '  -1.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0917'
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = '  +2.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.64'
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = '  +3.45%  '

$ws.Range("E35").Value = '  -0.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.64'
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = '  -5.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0351'
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = '  -0.54%  '

$ws.Range("E38").Value = '  -3.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.74'
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = '  +4.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.95'
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = '  +21.74%  '

$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.41'
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = '  +13.58%  '

$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.42'
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = '  +4.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.237'
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = '  +0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.26'
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = '  +19.45%  '

$ws.Range("E45").Value = '  +0.22%  '

$ws.Range("E46").Value = '  +0.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.67'
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = '  -1.15%  '

$ws.Range("E48").Value = '  -1.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.79'
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = '  +2.88%  '

$ws.Range("E50").Value = '  +1.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.452'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = '  +1.68%  '
